$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 3

$ws.Range("A1:A6").Select()
